$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 2374  # was 2367
$ws.Range("F3").Value = 559  # was 553
$ws.Range("F4").Value = 212  # was 211
$ws.Range("F5").Value = 367  # was 361
$ws.Range("F6").Value = 367  # was 361
$ws.Range("F7").Value = 604  # was 593
$ws.Range("F9").Value = 815  # was 804
$ws.Range("F10").Value = 543  # was 541
$ws.Range("F11").Value = 837  # was 833
$ws.Range("F13").Value = 102  # was 101
$ws.Range("F14").Value = 402  # was 403
$ws.Range("F16").Value = 1038  # was 1033
$ws.Range("F17").Value = 21739  # was 21526
$ws.Range("F18").Value = 944  # was 893
$ws.Range("F19").Value = 93  # was 83
$ws.Range("F20").Value = 282  # was 275
$ws.Range("F21").Value = 316  # was 308
$ws.Range("F22").Value = 182  # was 180
$ws.Range("F23").Value = 180  # was 171
$ws.Range("F25").Value = 21  # was 20
$ws.Range("F26").Value = 257  # was 253
$ws.Range("F28").Value = 370  # was 364
$ws.Range("F29").Value = 163  # was 162

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 5  # was 3
$ws.Range("F5").Value = 90  # was 89
$ws.Range("F6").Value = 211  # was 209
$ws.Range("F7").Value = 232  # was 230
$ws.Range("F8").Value = 3465  # was 3450
$ws.Range("F10").Value = 110  # was 107
$ws.Range("F13").Value = 29  # was 28
$ws.Range("F14").Value = 125  # was 124
$ws.Range("F16").Value = 3967  # was 3927

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 272  # was 270
$ws.Range("F3").Value = 123  # was 120
$ws.Range("F4").Value = 643  # was 637
$ws.Range("F5").Value = 218  # was 216

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 272  # was 270
$ws.Range("F3").Value = 123  # was 120
$ws.Range("F5").Value = 2374  # was 2367
$ws.Range("F6").Value = 643  # was 637
$ws.Range("F7").Value = 559  # was 553
$ws.Range("F8").Value = 212  # was 211
$ws.Range("F9").Value = 367  # was 361
$ws.Range("F10").Value = 367  # was 361
$ws.Range("F11").Value = 604  # was 593
$ws.Range("F14").Value = 5  # was 3
$ws.Range("F15").Value = 90  # was 89
$ws.Range("F16").Value = 211  # was 209
$ws.Range("F17").Value = 218  # was 216
$ws.Range("F18").Value = 815  # was 804
$ws.Range("F19").Value = 543  # was 541
$ws.Range("F20").Value = 837  # was 833
$ws.Range("F22").Value = 102  # was 101
$ws.Range("F23").Value = 402  # was 403
$ws.Range("F25").Value = 1038  # was 1033
$ws.Range("F26").Value = 21740  # was 21526
$ws.Range("F27").Value = 232  # was 230
$ws.Range("F28").Value = 3465  # was 3450
$ws.Range("F30").Value = 110  # was 107
$ws.Range("F32").Value = 944  # was 893
$ws.Range("F33").Value = 93  # was 83
$ws.Range("F34").Value = 282  # was 275
$ws.Range("F36").Value = 29  # was 28
$ws.Range("F37").Value = 316  # was 308
$ws.Range("F38").Value = 182  # was 180
$ws.Range("F39").Value = 180  # was 171
$ws.Range("F41").Value = 21  # was 20
$ws.Range("F42").Value = 125  # was 124
$ws.Range("F44").Value = 257  # was 253
$ws.Range("F46").Value = 370  # was 364
$ws.Range("F47").Value = 163  # was 162
$ws.Range("F48").Value = 3968  # was 3927

# Sheet 4: G26 changes from numeric 0 to text "已售罄" (sold out)
$ws = $wb.Worksheets.Item(4)
$ws.Range("G26").Value = "已售罄"
